$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0, 33.55332124156859, 39.681699737920887, 43.049412492974781, 101.20457063028111, 76.636095618834176, 50.411404645068792, 50.411404645068792),
    @(33.55332124156859, 0, 19.842385248134388, 84.666963306591285, 96.712188291840633, 62.917788541962778, 43.983000151641612, 43.983000151641612),
    @(39.681699737920887, 19.842385248134388, 0, 62.656054087204915, 103.42697049557917, 66.412242282720399, 52.881361297753301, 52.881361297753301),
    @(43.049412492974781, 84.666963306591285, 62.656054087204915, 0, 43.863701646047936, 51.136747775767461, 72.182304891540497, 72.182304891540497),
    @(101.20457063028111, 96.712188291840633, 103.42697049557917, 43.863701646047936, 0, 32.83843964005483, 64.30086530843441, 64.30086530843441),
    @(76.636095618834176, 62.917788541962778, 66.412242282720399, 51.136747775767461, 32.83843964005483, 0, 31.155806452025956, 31.155806452025956),
    @(50.411404645068792, 43.983000151641612, 52.881361297753301, 72.182304891540497, 64.30086530843441, 31.155806452025956, 0, 0),
    @(50.411404645068792, 43.983000151641612, 52.881361297753301, 72.182304891540497, 64.30086530843441, 31.155806452025956, 0, 0)
)

for ($r = 0; $r -lt 8; $r++) {
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $values[$r][$c]
    }
}
